$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the formulas in B27:B30 with a hardcoded (higher) value.
$ws.Range("B27:B30").Value = 10.337249999999999

# Update the active selection to match the edited cell.
$ws.Range("B27").Select()
